# ----------------------------------------------------------------------------
# asia_consolidated.xlsx regional data-processing pass:
#   1. Normalise the header row to lower_snake-ish labels and add units.
#   2. Shift the numeric measurement columns one slot to the left
#      (D <- old E, E <- old F) and populate the (now-empty) F column with
#      freshly computed "climate change (kg CO2 eq)" values.
#   3. Attach a data-dictionary cell comment to every header cell.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row relabeling -------------------------------------------------
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- 2. Data-dictionary comments on the header cells --------------------------
$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null

# --- 3. Column shift + newly computed "climate change" values (rows 2-62) -----

$ws.Range("D2").Value = 0.26148888
$ws.Range("E2").Value = 6.5711769
$ws.Range("F2").Value = 0.0000072910037
$ws.Range("D3").Value = 0.2458252266666667
$ws.Range("E3").Value = 6.3117561
$ws.Range("F3").Value = 0.0000068542596
$ws.Range("D4").Value = 0.2027121
$ws.Range("E4").Value = 5.5996635
$ws.Range("F4").Value = 0.0000056521513
$ws.Range("D5").Value = 0.2254119133333334
$ws.Range("E5").Value = 5.9748933
$ws.Range("F5").Value = 0.0000062850823
$ws.Range("D6").Value = 0.2106227666666667
$ws.Range("E6").Value = 5.730846
$ws.Range("F6").Value = 0.0000058727216
$ws.Range("D7").Value = 0.15633982
$ws.Range("E7").Value = 5.9187545
$ws.Range("F7").Value = 0.0000043591691
$ws.Range("D8").Value = 0.1655053733333333
$ws.Range("E8").Value = 4.9863687
$ws.Range("F8").Value = 0.0000046147289
$ws.Range("D9").Value = 0.1959247933333333
$ws.Range("E9").Value = 5.488086
$ws.Range("F9").Value = 0.0000054629031
$ws.Range("D10").Value = 0.25970618
$ws.Range("E10").Value = 6.5416111
$ws.Range("F10").Value = 0.0000072412974
$ws.Range("D11").Value = 0.21788696
$ws.Range("E11").Value = 6.2249096
$ws.Range("F11").Value = 0.0000060752667
$ws.Range("D12").Value = 0.1981833333333334
$ws.Range("E12").Value = 6.2476254
$ws.Range("F12").Value = 0.0000055258772
$ws.Range("D13").Value = 0.24230588
$ws.Range("E13").Value = 6.2539617
$ws.Range("F13").Value = 0.0000067561308
$ws.Range("D14").Value = 0.1605092933333334
$ws.Range("E14").Value = 6.0187442
$ws.Range("F14").Value = 0.000004475425
$ws.Range("D15").Value = 0.2191823466666667
$ws.Range("E15").Value = 5.8721818
$ws.Range("F15").Value = 0.0000061113854
$ws.Range("D16").Value = 0.2248807133333333
$ws.Range("E16").Value = 6.1115408
$ws.Range("F16").Value = 0.0000062702709
$ws.Range("D17").Value = 0.22132368
$ws.Range("E17").Value = 5.9079018
$ws.Range("F17").Value = 0.0000061710915
$ws.Range("D18").Value = 0.1235705066666667
$ws.Range("E18").Value = 4.2956038
$ws.Range("F18").Value = 0.0000034454736
$ws.Range("D19").Value = 0.1528896466666667
$ws.Range("E19").Value = 4.7792129
$ws.Range("F19").Value = 0.0000042629691
$ws.Range("D20").Value = 0.1943637333333333
$ws.Range("E20").Value = 6.2526771
$ws.Range("F20").Value = 0.0000054193766
$ws.Range("D21").Value = 0.1420520466666667
$ws.Range("E21").Value = 5.1634653
$ws.Range("F21").Value = 0.000003960788
$ws.Range("D22").Value = 0.18212972
$ws.Range("E22").Value = 6.2233437
$ws.Range("F22").Value = 0.0000050782598
$ws.Range("D23").Value = 0.2068272133333333
$ws.Range("E23").Value = 5.669548
$ws.Range("F23").Value = 0.0000057668915
$ws.Range("D24").Value = 0.04776863133333333
$ws.Range("E24").Value = 3.0464557
$ws.Range("F24").Value = 0.0000013319162
$ws.Range("D25").Value = 0.1789799533333333
$ws.Range("E25").Value = 5.2091253
$ws.Range("F25").Value = 0.000004990436
$ws.Range("D26").Value = 0.038210444
$ws.Range("E26").Value = 2.8880417
$ws.Range("F26").Value = 0.0000010654086
$ws.Range("D27").Value = 0.010111274
$ws.Range("E27").Value = 2.4246647
$ws.Range("F27").Value = 0.00000028192915
$ws.Range("D28").Value = 0.23257274
$ws.Range("E28").Value = 6.0932509
$ws.Range("F28").Value = 0.0000064847451
$ws.Range("D29").Value = 0.1337171266666667
$ws.Range("E29").Value = 4.4608343
$ws.Range("F29").Value = 0.0000037283882
$ws.Range("D30").Value = 0.03551259933333334
$ws.Range("E30").Value = 2.8387454
$ws.Range("F30").Value = 0.00000099018549
$ws.Range("D31").Value = 0.2029564333333333
$ws.Range("E31").Value = 5.6029334
$ws.Range("F31").Value = 0.0000056589639
$ws.Range("D32").Value = 0.1970088333333333
$ws.Range("E32").Value = 5.5059694
$ws.Range("F32").Value = 0.000005493129
$ws.Range("D33").Value = 0.3194690133333333
$ws.Range("E33").Value = 3.4160813
$ws.Range("F33").Value = 0.0000089076437
$ws.Range("D34").Value = 0.2765348133333334
$ws.Range("E34").Value = 3.2140221
$ws.Range("F34").Value = 0.0000077105244
$ws.Range("D35").Value = 0.002479979733333334
$ws.Range("E35").Value = 1.2396245
$ws.Range("F35").Value = 0.000000069148416
$ws.Range("D36").Value = 0.2740187666666667
$ws.Range("E36").Value = 3.1470098
$ws.Range("F36").Value = 0.0000076403703
$ws.Range("D37").Value = 0.34985066
$ws.Range("E37").Value = 3.7470553
$ws.Range("F37").Value = 0.0000097547647
$ws.Range("D38").Value = 0.003940904333333333
$ws.Range("E38").Value = 1.2396245
$ws.Range("F38").Value = 0.00000010988287
$ws.Range("D39").Value = 0.3459967333333334
$ws.Range("E39").Value = 3.737914
$ws.Range("F39").Value = 0.000009647307
$ws.Range("D40").Value = 0.002545318066666667
$ws.Range("E40").Value = 2.3796327
$ws.Range("F40").Value = 0.000000070970222
$ws.Range("D41").Value = 0.1937403933333333
$ws.Range("E41").Value = 2.5154459
$ws.Range("F41").Value = 0.0000054019963
$ws.Range("D42").Value = 0.003447503133333333
$ws.Range("E42").Value = 1.2473845
$ws.Range("F42").Value = 0.000000096125534
$ws.Range("D43").Value = 0.2079049066666667
$ws.Range("E43").Value = 2.8446502
$ws.Range("F43").Value = 0.0000057969405
$ws.Range("D44").Value = 0.3340773733333334
$ws.Range("E44").Value = 3.6653597
$ws.Range("F44").Value = 0.0000093149636
$ws.Range("D45").Value = 0.0024814414
$ws.Range("E45").Value = 1.2396245
$ws.Range("F45").Value = 0.00000006918917
$ws.Range("D46").Value = 0.002479248933333333
$ws.Range("E46").Value = 1.2396245
$ws.Range("F46").Value = 0.000000069128038
$ws.Range("D47").Value = 0.3496074733333334
$ws.Range("E47").Value = 3.7379459
$ws.Range("F47").Value = 0.0000097479841
$ws.Range("D48").Value = 0.1666445933333333
$ws.Range("E48").Value = 2.7229208
$ws.Range("F48").Value = 0.0000046464934
$ws.Range("D49").Value = 0.0026977714
$ws.Range("E49").Value = 1.2703211
$ws.Range("F49").Value = 0.000000075221026
$ws.Range("D50").Value = 0.3513559333333334
$ws.Range("E50").Value = 3.7505192
$ws.Range("F50").Value = 0.0000097967357
$ws.Range("D51").Value = 0.317259
$ws.Range("E51").Value = 3.5083763
$ws.Range("F51").Value = 0.0000088460228
$ws.Range("D52").Value = 0.28171948
$ws.Range("E52").Value = 3.4425417
$ws.Range("F52").Value = 0.0000078550866
$ws.Range("D53").Value = 0.0025019694
$ws.Range("E53").Value = 1.2400211
$ws.Range("F53").Value = 0.000000069761546
$ws.Range("D54").Value = 0.002479248933333333
$ws.Range("E54").Value = 1.2396245
$ws.Range("F54").Value = 0.000000069128038
$ws.Range("D55").Value = 0.002500589066666667
$ws.Range("E55").Value = 1.2396245
$ws.Range("F55").Value = 0.000000069723059
$ws.Range("D56").Value = 0.002479248933333333
$ws.Range("E56").Value = 1.2396245
$ws.Range("F56").Value = 0.000000069128038
$ws.Range("D57").Value = 0.3177190333333333
$ws.Range("E57").Value = 3.5123024
$ws.Range("F57").Value = 0.0000088588497
$ws.Range("D58").Value = 0.21942052
$ws.Range("E58").Value = 2.6768029
$ws.Range("F58").Value = 0.0000061180263
$ws.Range("D59").Value = 0.2726714666666667
$ws.Range("E59").Value = 3.2398612
$ws.Range("F59").Value = 0.000007602804
$ws.Range("D60").Value = 0.1860564533333333
$ws.Range("E60").Value = 2.7124628
$ws.Range("F60").Value = 0.0000051877475
$ws.Range("D61").Value = 0.002479248933333333
$ws.Range("E61").Value = 1.2396245
$ws.Range("F61").Value = 0.000000069128038
$ws.Range("D62").Value = 0.2073880133333333
$ws.Range("E62").Value = 3.0687868
$ws.Range("F62").Value = 0.0000057825281
